$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2824.5
$ws.Range("I40").Value = 2699.5
$ws.Range("K40").Value = 2699.5
$ws.Range("M40").Value = -2524.5
$ws.Range("H62").Value = 1699.3334
$ws.Range("J62").Value = 1699
$ws.Range("L62").Value = 1699
$ws.Range("N62").Value = -2947
$ws.Range("H65").Value = 1699.3334
$ws.Range("J65").Value = 1699
$ws.Range("L65").Value = 8495
$ws.Range("N65").Value = -14735
$ws.Range("H70").Value = 1750
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1750
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5250
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -5790
$ws.Range("H73").Value = 1750
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1750
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5250
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -7122
$ws.Range("H112").Value = 6995
$ws.Range("J112").Value = 6995
$ws.Range("L112").Value = 20985
$ws.Range("N112").Value = -23201
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -8508
$ws.Range("H116").Value = 8000
$ws.Range("I116").Value = 8000
$ws.Range("K116").Value = 8000
$ws.Range("M116").Value = -4558
$ws.Range("H137").Value = 2791.4285
$ws.Range("I137").Value = 2882.6667
$ws.Range("J137").Value = 2244
$ws.Range("K137").Value = 8648.000100000001
$ws.Range("L137").Value = 6732
$ws.Range("M137").Value = -6098.000100000001
$ws.Range("N137").Value = -11832
$ws.Range("H138").Value = 2775

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H124").Value = 33660.75
$ws.Range("J124").Value = 33660.75
$ws.Range("L124").Value = 33660.75
$ws.Range("N124").Value = -43480.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 40
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 40
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 73
$ws.Range("N6").Value = ""
$ws.Range("H16").Value = 1316
$ws.Range("I16").Value = 1270
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1270
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -983
$ws.Range("N16").Value = -2074
$ws.Range("H28").Value = 51599.668
$ws.Range("J28").Value = 51599.668
$ws.Range("L28").Value = 51599.668
$ws.Range("N28").Value = -52089.668
$ws.Range("H99").Value = 50000
$ws.Range("I99").Value = 50000
$ws.Range("K99").Value = 50000
$ws.Range("M99").Value = -48502
$ws.Range("H113").Value = 1316
$ws.Range("I113").Value = 1270
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1270
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 900
$ws.Range("N113").Value = -5840
$ws.Range("H126").Value = 50000
$ws.Range("I126").Value = 50000
$ws.Range("K126").Value = 150000
$ws.Range("M126").Value = -147530

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22.5
$ws.Range("J2").Value = 20
$ws.Range("L2").Value = 120
$ws.Range("N2").Value = -346
$ws.Range("H7").Value = 79.40000000000001
$ws.Range("J7").Value = 155.5
$ws.Range("L7").Value = 466.5
$ws.Range("N7").Value = -690.5
$ws.Range("H68").Value = 2497.5
$ws.Range("J68").Value = 2497.5
$ws.Range("L68").Value = 7492.5
$ws.Range("N68").Value = -9114.5
$ws.Range("H71").Value = 2497.5
$ws.Range("J71").Value = 2497.5
$ws.Range("L71").Value = 22477.5
$ws.Range("N71").Value = -30589.5
$ws.Range("H80").Value = 13725
$ws.Range("J80").Value = 13725
$ws.Range("L80").Value = 41175
$ws.Range("N80").Value = -43047
$ws.Range("H83").Value = 13725
$ws.Range("J83").Value = 13725
$ws.Range("L83").Value = 123525
$ws.Range("N83").Value = -132885
$ws.Range("H109").Value = 1933
$ws.Range("J109").Value = 2000
$ws.Range("L109").Value = 6000
$ws.Range("N109").Value = -8080
$ws.Range("H117").Value = 2125
$ws.Range("J117").Value = 2125
$ws.Range("L117").Value = 6375
$ws.Range("N117").Value = -13259
$ws.Range("H132").Value = 1999
$ws.Range("I132").Value = 1999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17991
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15461
$ws.Range("N132").Value = ""
$ws.Range("H141").Value = 2866.6667
$ws.Range("J141").Value = 3000
$ws.Range("L141").Value = 9000
$ws.Range("N141").Value = -19360

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -29984
$ws.Range("H135").Value = 91998.664
$ws.Range("J135").Value = 91998.664
$ws.Range("L135").Value = 91998.664
$ws.Range("N135").Value = -102138.664

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5350.5
$ws.Range("I40").Value = 5350.5
$ws.Range("K40").Value = 5350.5
$ws.Range("M40").Value = -5214.5
$ws.Range("H82").Value = 1666.6666
$ws.Range("I82").Value = 1666.6666
$ws.Range("K82").Value = 1666.6666
$ws.Range("M82").Value = -1305.6666
$ws.Range("H85").Value = 1666.6666
$ws.Range("I85").Value = 1666.6666
$ws.Range("K85").Value = 1666.6666
$ws.Range("M85").Value = -418.6666

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 952.7
$ws.Range("J113").Value = 1950
$ws.Range("L113").Value = 5850
$ws.Range("N113").Value = -10190
$ws.Range("H136").Value = 1927.6428
$ws.Range("I136").Value = 1927.6428
$ws.Range("K136").Value = 5782.928400000001
$ws.Range("M136").Value = -3232.928400000001
